$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": column N (27-jun) for hours 00-01 .. 23-24 (rows 2..25) ---
# These readings are not available yet -> replace the numeric values with "-"
$wsSpot = $wb.Worksheets.Item("Prix Spot")
for ($r = 2; $r -le 25; $r++) {
    $wsSpot.Cells.Item($r, 14).Value = "-"
}

# --- Sheet "CO2": insert two missing daily quotes (2025-06-21 and 2025-06-22) ---
# before: ... 2025-06-20 (row6), 2025-06-23 (row7), 2025-06-25 (row8)
# after : ... 2025-06-20 (row6), 2025-06-21 (row7), 2025-06-22 (row8), 2025-06-23 (row9), 2025-06-25 (row10)
$wsCo2 = $wb.Worksheets.Item("CO2")

# push the existing 2025-06-23 / 2025-06-25 rows down by two
$wsCo2.Rows.Item(7).Insert()
$wsCo2.Rows.Item(7).Insert()

# the "Date" column stores plain text dates (not Excel date serials), so force
# text entry with a leading apostrophe and strip the resulting cell format so
# the new cells end up looking just like their neighbours
$a7 = $wsCo2.Cells.Item(7, 1)
$a7.Value = "'2025-06-21"
$a7.ClearFormats()
$wsCo2.Cells.Item(7, 2).Value = 72.2

$a8 = $wsCo2.Cells.Item(8, 1)
$a8.Value = "'2025-06-22"
$a8.ClearFormats()
$wsCo2.Cells.Item(8, 2).Value = 72.2
